# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure a reusable "text" number format cell style trick is available.
# For numeric-looking D-column values we must force text storage so that
# values such as "0.530" keep their trailing zero instead of becoming 0.53.

# Row 2
$ws.Range("D2").Value = "64.391.04"
$ws.Range("E2").Value = "  +0.10%  "

# Row 3
$ws.Range("D3").Value = "3.135.80"
$ws.Range("E3").Value = "  -0.66%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.39%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.68%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").Value = "3.131.50"
$ws.Range("E8").Value = "  -0.69%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = "Normal"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.65%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.35%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.470"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.49%  "

# Row 13
$ws.Range("E13").Value = "  +1.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "

# Row 15
$ws.Range("D15").Value = "3.654.65"
$ws.Range("E15").Value = "  -0.66%  "

# Row 16
$ws.Range("E16").Value = "  +2.33%  "

# Row 17
$ws.Range("D17").Value = "64.358.91"
$ws.Range("E17").Value = "  +0.02%  "

# Row 18
$ws.Range("D18").Value = "3.141.65"
$ws.Range("E18").Value = "  -0.50%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.08%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.27%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.47%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.717"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.02%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.79"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.19%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.57%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.93%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.59%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.85%  "

# Row 30
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.06%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.115"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "

# Row 32
$ws.Range("E32").Value = "  -0.12%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.26%  "

# Row 34
$ws.Range("E34").Value = "  -3.39%  "

# Row 35
$ws.Range("E35").Value = "  +0.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.97"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.20%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.56%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0742"
$ws.Range("E38").Value = "  +1.59%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "450.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.37%  "

# Row 40
$ws.Range("E40").Value = "  +1.66%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.118"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.28%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.32"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.26%  "

# Row 44
$ws.Range("D44").Value = "2.875.12"
$ws.Range("E44").Value = "  +0.63%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.262"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.70%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.65%  "

# Row 50
$ws.Range("E50").Value = "  -0.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.97%  "
